# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for each coin row on the active sheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2 = @{ D='26.427.70'; E='  -1.07%  ' }
    3 = @{ D='1.624.65'; E='  -0.69%  ' }
    4 = @{ E='  +0.08%  ' }
    5 = @{ D='212.74'; E='  -0.29%  ' }
    6 = @{ D='0.498'; E='  +1.66%  ' }
    7 = @{ E='  +0.03%  ' }
    8 = @{ E='  -0.98%  ' }
    9 = @{ E='  +0.34%  ' }
    10 = @{ D='18.97'; E='  -0.68%  ' }
    11 = @{ D='0.0844'; E='  +1.13%  ' }
    12 = @{ D='1.850.50' }
    13 = @{ D='4.15'; E='  +2.23%  ' }
    14 = @{ D='1.614.07'; E='  -1.72%  ' }
    15 = @{ E='  -0.49%  ' }
    16 = @{ D='64.37'; E='  +1.92%  ' }
    17 = @{ D='26.438.89'; E='  -1.00%  ' }
    18 = @{ E='  +0.19%  ' }
    19 = @{ D='216.11'; E='  +3.81%  ' }
    20 = @{ E='  +0.22%  ' }
    21 = @{ E='  -0.53%  ' }
    22 = @{ D='6.23'; E='  +2.03%  ' }
    23 = @{ D='9.29'; E='  -0.92%  ' }
    24 = @{ D='1.99'; E='  +4.79%  ' }
    25 = @{ D='147.73'; E='  +0.99%  ' }
    26 = @{ E='  +0.13%  ' }
    27 = @{ E='  -0.97%  ' }
    28 = @{ D='6.85'; E='  +2.56%  ' }
    29 = @{ E='  +1.16%  ' }
    30 = @{ D='0.0510'; E='  -1.95%  ' }
    31 = @{ E='  -1.21%  ' }
    32 = @{ D='3.30'; E='  +2.22%  ' }
    33 = @{ E='  -0.48%  ' }
    34 = @{ E='  -0.82%  ' }
    35 = @{ D='1.218.54'; E='  +4.49%  ' }
    36 = @{ E='  -1.86%  ' }
    37 = @{ D='0.0172'; E='  +3.24%  ' }
    38 = @{ E='  +0.06%  ' }
    39 = @{ D='0.794'; E='  -1.92%  ' }
    40 = @{ D='0.502'; E='  +0.16%  ' }
    41 = @{ E='  -3.00%  ' }
    42 = @{ E='  -0.41%  ' }
    43 = @{ E='  +0.12%  ' }
    44 = @{ D='1.761.87'; E='  -0.86%  ' }
    45 = @{ D='92.80'; E='  +0.28%  ' }
    46 = @{ E='  +1.82%  ' }
    47 = @{ D='54.69'; E='  +0.14%  ' }
    48 = @{ E='  -2.13%  ' }
    49 = @{ D='0.0509'; E='  -0.55%  ' }
    50 = @{ D='7.56'; E='  -1.07%  ' }
    51 = @{ D='0.406'; E='  -0.96%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey('D')) {
        $cell = $ws.Range("D$row")
        # Force text storage (matches the workbook's existing inline-string
        # cells) so values like "212.74" or "18.97" aren't reinterpreted as
        # numbers; then drop back to the default style so no stray
        # number-format is left applied to the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $vals['D']
        $cell.Style = "Normal"
    }

    if ($vals.ContainsKey('E')) {
        $ws.Range("E$row").Value = $vals['E']
    }
}
